$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efna2"
$ws.Cells.Item(2,3).Value = "Epha5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.7317105
$ws.Cells.Item(2,8).Value = 3.463421
$ws.Cells.Item(2,9).Value = 0.2580726189025833
$ws.Cells.Item(2,10).Value = 0.2162491041706965
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.5
$ws.Cells.Item(2,13).Value = 0.006466
$ws.Cells.Item(2,14).Value = 0.012932
$ws.Cells.Item(2,15).Value = 0.008493427970384656
$ws.Cells.Item(2,16).Value = 0.008493427970384656
$ws.Cells.Item(2,17).Value = 0.011197240093
$ws.Cells.Item(2,18).Value = 0.04478896037200001
$ws.Cells.Item(2,19).Value = 0.002191921199777621
$ws.Cells.Item(2,20).Value = 0.001836696189934019

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efna2"
$ws.Cells.Item(3,3).Value = "Epha5"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.7317105
$ws.Cells.Item(3,8).Value = 3.463421
$ws.Cells.Item(3,9).Value = 0.2580726189025833
$ws.Cells.Item(3,10).Value = 0.2162491041706965
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.7548284999999999
$ws.Cells.Item(3,14).Value = 1.509657
$ws.Cells.Item(3,15).Value = 0.9915065720296153
$ws.Cells.Item(3,16).Value = 0.9915065720296153
$ws.Cells.Item(3,17).Value = 1.30714443914925
$ws.Cells.Item(3,18).Value = 5.228577756597
$ws.Cells.Item(3,19).Value = 0.2558806977028056
$ws.Cells.Item(3,20).Value = 0.2144124079807625

$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Efna2"
$ws.Cells.Item(4,3).Value = "Epha5"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.309486666666667
$ws.Cells.Item(4,8).Value = 6.92846
$ws.Cells.Item(4,9).Value = 0.3441772007430018
$ws.Cells.Item(4,10).Value = 0.4325992330364988
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.5
$ws.Cells.Item(4,13).Value = 0.006466
$ws.Cells.Item(4,14).Value = 0.012932
$ws.Cells.Item(4,15).Value = 0.008493427970384656
$ws.Cells.Item(4,16).Value = 0.008493427970384656
$ws.Cells.Item(4,17).Value = 0.01493314078666667
$ws.Cells.Item(4,18).Value = 0.08959884472000001
$ws.Cells.Item(4,19).Value = 0.002923244263559307
$ws.Cells.Item(4,20).Value = 0.003674250425839149

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Efna2"
$ws.Cells.Item(5,3).Value = "Epha5"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.309486666666667
$ws.Cells.Item(5,8).Value = 6.92846
$ws.Cells.Item(5,9).Value = 0.3441772007430018
$ws.Cells.Item(5,10).Value = 0.4325992330364988
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.7548284999999999
$ws.Cells.Item(5,14).Value = 1.509657
$ws.Cells.Item(5,15).Value = 0.9915065720296153
$ws.Cells.Item(5,16).Value = 0.9915065720296153
$ws.Cells.Item(5,17).Value = 1.74326635637
$ws.Cells.Item(5,18).Value = 10.45959813822
$ws.Cells.Item(5,19).Value = 0.3412539564794425
$ws.Cells.Item(5,20).Value = 0.4289249826106597

$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Efna2"
$ws.Cells.Item(6,3).Value = "Epha5"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.382906
$ws.Cells.Item(6,8).Value = 4.765812
$ws.Cells.Item(6,9).Value = 0.3551187060531648
$ws.Cells.Item(6,10).Value = 0.2975678023682237
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.5
$ws.Cells.Item(6,13).Value = 0.006466
$ws.Cells.Item(6,14).Value = 0.012932
$ws.Cells.Item(6,15).Value = 0.008493427970384656
$ws.Cells.Item(6,16).Value = 0.008493427970384656
$ws.Cells.Item(6,17).Value = 0.015407870196
$ws.Cells.Item(6,18).Value = 0.06163148078400001
$ws.Cells.Item(6,19).Value = 0.003016175150798757
$ws.Cells.Item(6,20).Value = 0.002527370695720164

$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Efna2"
$ws.Cells.Item(7,3).Value = "Epha5"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.382906
$ws.Cells.Item(7,8).Value = 4.765812
$ws.Cells.Item(7,9).Value = 0.3551187060531648
$ws.Cells.Item(7,10).Value = 0.2975678023682237
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7548284999999999
$ws.Cells.Item(7,14).Value = 1.509657
$ws.Cells.Item(7,15).Value = 0.9915065720296153
$ws.Cells.Item(7,16).Value = 0.9915065720296153
$ws.Cells.Item(7,17).Value = 1.798685361621
$ws.Cells.Item(7,18).Value = 7.194741446484
$ws.Cells.Item(7,19).Value = 0.352102530902366
$ws.Cells.Item(7,20).Value = 0.2950404316725035

$ws.Cells.Item(8,1).Value = "Neutrophils"
$ws.Cells.Item(8,2).Value = "Efna2"
$ws.Cells.Item(8,3).Value = "Epha5"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.05833266666666666
$ws.Cells.Item(8,8).Value = 0.174998
$ws.Cells.Item(8,9).Value = 0.008693175940342274
$ws.Cells.Item(8,10).Value = 0.01092652632517489
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.5
$ws.Cells.Item(8,13).Value = 0.006466
$ws.Cells.Item(8,14).Value = 0.012932
$ws.Cells.Item(8,15).Value = 0.008493427970384656
$ws.Cells.Item(8,16).Value = 0.008493427970384656
$ws.Cells.Item(8,17).Value = 0.0003771790226666667
$ws.Cells.Item(8,18).Value = 0.002263074136
$ws.Cells.Item(8,19).Value = 0.00007383486368317799
$ws.Cells.Item(8,20).Value = 0.000092803664309384667

$ws.Cells.Item(9,1).Value = "Neutrophils"
$ws.Cells.Item(9,2).Value = "Efna2"
$ws.Cells.Item(9,3).Value = "Epha5"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.05833266666666666
$ws.Cells.Item(9,8).Value = 0.174998
$ws.Cells.Item(9,9).Value = 0.008693175940342274
$ws.Cells.Item(9,10).Value = 0.01092652632517489
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.7548284999999999
$ws.Cells.Item(9,14).Value = 1.509657
$ws.Cells.Item(9,15).Value = 0.9915065720296153
$ws.Cells.Item(9,16).Value = 0.9915065720296153
$ws.Cells.Item(9,17).Value = 0.04403115928099999
$ws.Cells.Item(9,18).Value = 0.2641869556859999
$ws.Cells.Item(9,19).Value = 0.008619341076659095
$ws.Cells.Item(9,20).Value = 0.01083372266086551

$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Efna2"
$ws.Cells.Item(10,3).Value = "Epha5"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.2277316666666667
$ws.Cells.Item(10,8).Value = 0.683195
$ws.Cells.Item(10,9).Value = 0.03393829836090777
$ws.Cells.Item(10,10).Value = 0.04265733409940604
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.5
$ws.Cells.Item(10,13).Value = 0.006466
$ws.Cells.Item(10,14).Value = 0.012932
$ws.Cells.Item(10,15).Value = 0.008493427970384656
$ws.Cells.Item(10,16).Value = 0.008493427970384656
$ws.Cells.Item(10,17).Value = 0.001472512956666667
$ws.Cells.Item(10,18).Value = 0.00883507774
$ws.Cells.Item(10,19).Value = 0.0002882524925657938
$ws.Cells.Item(10,20).Value = 0.0003623069945819385

$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Efna2"
$ws.Cells.Item(11,3).Value = "Epha5"
$ws.Cells.Item(11,4).Value = "MuSCs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.2277316666666667
$ws.Cells.Item(11,8).Value = 0.683195
$ws.Cells.Item(11,9).Value = 0.03393829836090777
$ws.Cells.Item(11,10).Value = 0.04265733409940604
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.7548284999999999
$ws.Cells.Item(11,14).Value = 1.509657
$ws.Cells.Item(11,15).Value = 0.9915065720296153
$ws.Cells.Item(11,16).Value = 0.9915065720296153
$ws.Cells.Item(11,17).Value = 0.1718983523525
$ws.Cells.Item(11,18).Value = 1.031390114115
$ws.Cells.Item(11,19).Value = 0.03365004586834198
$ws.Cells.Item(11,20).Value = 0.04229502710482411
